# Add a new weekly price record as row 7, pushing the existing rows
# (old rows 7-33) down by one (new rows 8-34).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 7 - this shifts rows 7:33 down to 8:34
# and extends the used range / dimension accordingly.
$ws.Rows.Item(7).Insert()

# Populate the new row 7 with the new weekly record.
$ws.Cells.Item(7, 1).Value = 7
$ws.Cells.Item(7, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(7, 3).Value = "Ñuble"
$ws.Cells.Item(7, 4).Value = 44532
$ws.Cells.Item(7, 5).Value = 16
$ws.Cells.Item(7, 6).Value = 100112022
$ws.Cells.Item(7, 7).Value = "Arveja Verde"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 60
$ws.Cells.Item(7, 11).Value = 13000
$ws.Cells.Item(7, 12).Value = 14000
$ws.Cells.Item(7, 13).Value = 13500
$ws.Cells.Item(7, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(7, 15).Value = "Región del Maule"
$ws.Cells.Item(7, 16).Value = 540
$ws.Cells.Item(7, 17).Value = 25
$ws.Cells.Item(7, 18).Value = "Hortaliza"
